# Add two new columns (I: "I0", J: "IF") to the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row -- copy the existing header formatting (bold, centered,
# bordered) from H1 onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-21: column I ("I0") and column J ("IF") values.
# For every row J mirrors column H ("IP"), and I is 1 -- except row 20,
# which carries distinct values (I20=4, J20=6).
$iValues = @{
    2  = 1;  3  = 1;  4  = 1;  5  = 1;  6  = 1;  7  = 1;  8  = 1;  9  = 1;
    10 = 1;  11 = 1;  12 = 1;  13 = 1;  14 = 1;  15 = 1;  16 = 1;  17 = 1;
    18 = 1;  19 = 1;  20 = 4;  21 = 1
}
$jValues = @{
    2  = 5;  3  = 5;  4  = 6;  5  = 6;  6  = 5;  7  = 4;  8  = 6;  9  = 5;
    10 = 4;  11 = 2;  12 = 5;  13 = 7;  14 = 6;  15 = 3;  16 = 4;  17 = 5;
    18 = 5;  19 = 4;  20 = 6;  21 = 2
}

for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
